$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume update scraped on Tue Sep 17 03:00:14 UTC 2024

$ws.Range("D2").Value = "58.120.78"
$ws.Range("E2").Value = "  -0.17%  "
$ws.Range("D3").Value = "2.286.18"
$ws.Range("E3").Value = "  +1.26%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "534.04"
$ws.Range("E5").Value = "  -1.40%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.75"
$ws.Range("E6").Value = "  +0.99%  "
$ws.Range("E7").Value = "  -0.22%  "
$ws.Range("E8").Value = "  +4.27%  "
$ws.Range("D9").Value = "2.283.90"
$ws.Range("E9").Value = "  +1.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0994"
$ws.Range("E10").Value = "  -0.77%  "
$ws.Range("E11").Value = "  +0.92%  "
$ws.Range("E12").Value = "  +0.98%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.331"
$ws.Range("E13").Value = "  +0.00%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.40"
$ws.Range("E14").Value = "  +0.15%  "
$ws.Range("D15").Value = "2.695.18"
$ws.Range("E15").Value = "  +0.81%  "
$ws.Range("D16").Value = "58.056.52"
$ws.Range("E16").Value = "  -0.28%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000131"
$ws.Range("E17").Value = "  +0.32%  "
$ws.Range("D18").Value = "2.272.98"
$ws.Range("E18").Value = "  -0.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.48"
$ws.Range("E19").Value = "  -0.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.17"
$ws.Range("E20").Value = "  -1.70%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "313.20"
$ws.Range("E21").Value = "  +0.47%  "
$ws.Range("E22").Value = "  +0.77%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.65"
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.95"
$ws.Range("E27").Value = "  -0.92%  "
$ws.Range("E28").Value = "  -2.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "170.93"
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("E30").Value = "  -1.59%  "
$ws.Range("D31").Value = "0.0₃0721"
$ws.Range("E31").Value = "  +1.65%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.75"
$ws.Range("E32").Value = "  -0.22%  "
$ws.Range("E33").Value = "  -0.94%  "
$ws.Range("E34").Value = "  +0.27%  "
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.78"
$ws.Range("E36").Value = "  +1.00%  "
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("E38").Value = "  -0.30%  "
$ws.Range("E39").Value = "  -0.19%  "
$ws.Range("E40").Value = "  -0.77%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "140.60"
$ws.Range("E41").Value = "  +1.51%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "286.64"
$ws.Range("E42").Value = "  -1.95%  "
$ws.Range("E43").Value = "  +0.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0952"
$ws.Range("E44").Value = "  +1.37%  "
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.552"
$ws.Range("E46").Value = "  +1.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.90"
$ws.Range("E47").Value = "  -0.66%  "
$ws.Range("E48").Value = "  -0.83%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.95"
$ws.Range("E49").Value = "  -0.48%  "
$ws.Range("E50").Value = "  -0.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.51"
$ws.Range("E51").Value = "  +2.55%  "
